$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Meeting time: "9:00 -10:00 am" -> "9:00 -11:00 am"
#    ("10:00" becomes "11:00" - found via Find, then only the tail
#    "0:00" is rewritten to "1:00", leaving the leading "1" alone.)
# ------------------------------------------------------------------
$timeRng = $d.Content
$timeRng.Find.Execute("10:00", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$timeTail = $d.Range($timeRng.Start + 1, $timeRng.End)
$timeTail.Text = "1:00"

# ------------------------------------------------------------------
# 2) Add three new agenda bullets right before "New work assignment."
#    (same ListParagraph / numId 17 bullet list, Times New Roman
#    (Body CS) complex-script font, 12pt text)
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("New work assignment.", $false, $false, $false, `
                       $false, $false, $true, 1, $false, "", 0)
$assignmentIndex = $findRng.Paragraphs(1).Index

$newBullets = @(
    "Created JiangLab organization in MongoDB.",
    "Created a new project and database in MongoDB",
    "Updated the technote regarding MongoDB."
)

foreach ($bulletText in $newBullets) {
    $anchorPara = $d.Paragraphs($assignmentIndex)
    $insPoint = $anchorPara.Range
    $insPoint.Collapse(1)
    $insPoint.InsertParagraphBefore()

    $newPara = $d.Paragraphs($assignmentIndex)
    $newPara.Range.Text = $bulletText
    $newPara.Range.Font.NameBi = "Times New Roman (Body CS)"

    $assignmentIndex = $assignmentIndex + 1
}

# ------------------------------------------------------------------
# 3) Screenshots: mark the picture runs NoProof (adds <w:noProof/>)
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $d.InlineShapes($i).Range.NoProofing = $true
}

# ------------------------------------------------------------------
# 4) Drop the stray lastRenderedPageBreak that used to sit in front
#    of "Shap" (SHAP) - round-trip the text so the marker is cleared.
# ------------------------------------------------------------------
$shapRng = $d.Content
$shapRng.Find.Execute("Shap", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$shapRng.Text = "Shape"
$shapRng2 = $d.Content
$shapRng2.Find.Execute("Shape", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "Shap", 2)

Write-Host "Edit complete"
